# ORG: Lógica de controle trabalho 1 finalizada
#
# Adds a final verification row (row 14) below the truth table, writing
# "OK" under the first group of control signals (D:H) and "ok" under the
# second group (I:O), then leaves the selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data just below the existing table body (row 13 was the last).
$ws.Range("D14:H14").Value = "OK"
$ws.Range("I14:O14").Value = "ok"

# Leave the selection on the cell the author ended up on.
$null = $ws.Range("M23").Select()
